$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Param_Acq_MN")
$ws.Name = "MN_CQ_Gamma_cam"
